$d = $word.ActiveDocument

# 1. Merge the "arrival times" sentence back into a single run (remove the
#    proofErr-wrapped "0," split). The visible text doesn't change, just
#    the run structure simplifies, so a Find/Replace with identical text
#    achieves the run merge as a side effect of Word's replace.
$d.Content.Find.Execute(
    "Consider 5 jobs, A through E, with runtimes 3, 5, 2, 2, 2 and arrival times 0, 0, 5, 5, 5 respectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Consider 5 jobs, A through E, with runtimes 3, 5, 2, 2, 2 and arrival times 0, 0, 5, 5, 5 respectively.",
    2) | Out-Null

# 2. "using only one of the programming languages (C, C++, Java, Python)" -> "using Python"
$d.Content.Find.Execute(
    "using only one of the programming languages (C, C++, Java, Python)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "using Python",
    2) | Out-Null

# 3. Move the lastRenderedPageBreak marker: drop it from before "Functionality, 25 points"
#    and merge "2. " + "READMe" + " file (5 points): " into one run that now carries
#    the page-break marker.
$d.Content.Find.Execute(
    "Functionality, 25 points",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Functionality, 25 points",
    2) | Out-Null

$d.Content.Find.Execute(
    "2. READMe file (5 points): ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "2. READMe file (5 points): ",
    2) | Out-Null

# 4. Merge "solution as a single zipped file (." + "rar" + ", .zip) " into one run.
$d.Content.Find.Execute(
    "solution as a single zipped file (.rar, .zip) ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "solution as a single zipped file (.rar, .zip) ",
    2) | Out-Null

# 5. Merge the file-extension list run and the trailing README/Design doc run.
$d.Content.Find.Execute(
    "with the appropriate file extension like .java, .py, .c, .c++",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "with the appropriate file extension like .java, .py, .c, .c++",
    2) | Out-Null

$d.Content.Find.Execute(
    "), READMe file / Design document and pdf document. Non-adherence to instructions would result in losing points.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "), READMe file / Design document and pdf document. Non-adherence to instructions would result in losing points.",
    2) | Out-Null


